$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1750183612496/A5MTrC4t3W-camera-video-segment-driver-1750183614996.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSF4KUIPEM%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T010005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDWqpT3vXbTcJEykZGM%2F7X%2B8cCmNzx24gLc6N7vIEBARAIhALOJOs9Gn6j%2BXidxNTNrq7Bxn%2F5siLlT74MJ6mGD15mIKuYDCJj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgwwGDdKhFzPvUQY3wAqugPXUgcgeUu4qcQIkG3N7shPvpQW7j5wHUCsgW7V53Olq6t0LfSqJfNHT26wH14Bq12VcntUBVt5E2qrV2PBzUgFLpWipERqDdCJeoofeWt1XFjaXBLNYHgUhXAxcTvrMJPf5C5cycHyoE3KORYyzBwlCzNCwbexzk%2BM0iWoauNsOfRGqc9aElrjUtEFZcWh%2Fd5EWgsWXH3SyJQa0C%2F5GLaXRiQgxaS20fMxuwiLnKuD7SJilwLCLmHe7msZKNB2W3n4JXM3jTwbpEzya1PDZhxlXOvuwA%2BkMuZGC9i7PgdlzplqTSkk8sXzVFRtpZA4L1NCItMSBy8AzE5vTnhZJi6wpfOtswcg11UuBe2jdtSWO7qjXTWBHztMrWTVYdk8ghl9dE5vVgPdeeEAkW%2BI5OsOccb5YikxLxgZV1PWFhvdNjl%2FLtaIDTx%2FrG%2FlP9S%2Fwg9kBS1CF%2BzpXfQyPUYtvH8vpFNulkOMACWSQMlHDpAU%2B0sTJifj2plwRRLJYqsx96E8FP2U0ihLZpQC4up%2F1KHOL3EvCVQxbMCma5hOGQhWomUr%2BvfBhCCXgL%2F025TO5cBe4OVu9WkJZfrJMLD8zMIGOqQB7pW6hWXA2tomga2LqDXXQey05yP81WXFX1Yep26oiPmD5hCVc1H1A80szBa4GjLlmAYvMOMvJ3hcZHpy2vIj9KO70quMcwovPnCeBKaoO5DROHvq9oO3NRPMetMR9u5qmXrlQubOcL04T5Ez6O6k8QjzZpk2QvFn%2BgwqhYdol%2B%2FF%2F2UiFXBWwtFoAo5G4hY7qGEcLlncpi1PXUTqOuq6DgrvsyQ%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2009%3A00%3A05%20GMT&X-Amz-Signature=333b898b38510156ad8380c29b5a41f3a4c99aeb8f89395a22bb1f1a70e4e780"
$ws.Range("L2").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1750183612496/TAdN8NokD3-camera-video-segment-1750183614996.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSF4KUIPEM%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T010005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDWqpT3vXbTcJEykZGM%2F7X%2B8cCmNzx24gLc6N7vIEBARAIhALOJOs9Gn6j%2BXidxNTNrq7Bxn%2F5siLlT74MJ6mGD15mIKuYDCJj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgwwGDdKhFzPvUQY3wAqugPXUgcgeUu4qcQIkG3N7shPvpQW7j5wHUCsgW7V53Olq6t0LfSqJfNHT26wH14Bq12VcntUBVt5E2qrV2PBzUgFLpWipERqDdCJeoofeWt1XFjaXBLNYHgUhXAxcTvrMJPf5C5cycHyoE3KORYyzBwlCzNCwbexzk%2BM0iWoauNsOfRGqc9aElrjUtEFZcWh%2Fd5EWgsWXH3SyJQa0C%2F5GLaXRiQgxaS20fMxuwiLnKuD7SJilwLCLmHe7msZKNB2W3n4JXM3jTwbpEzya1PDZhxlXOvuwA%2BkMuZGC9i7PgdlzplqTSkk8sXzVFRtpZA4L1NCItMSBy8AzE5vTnhZJi6wpfOtswcg11UuBe2jdtSWO7qjXTWBHztMrWTVYdk8ghl9dE5vVgPdeeEAkW%2BI5OsOccb5YikxLxgZV1PWFhvdNjl%2FLtaIDTx%2FrG%2FlP9S%2Fwg9kBS1CF%2BzpXfQyPUYtvH8vpFNulkOMACWSQMlHDpAU%2B0sTJifj2plwRRLJYqsx96E8FP2U0ihLZpQC4up%2F1KHOL3EvCVQxbMCma5hOGQhWomUr%2BvfBhCCXgL%2F025TO5cBe4OVu9WkJZfrJMLD8zMIGOqQB7pW6hWXA2tomga2LqDXXQey05yP81WXFX1Yep26oiPmD5hCVc1H1A80szBa4GjLlmAYvMOMvJ3hcZHpy2vIj9KO70quMcwovPnCeBKaoO5DROHvq9oO3NRPMetMR9u5qmXrlQubOcL04T5Ez6O6k8QjzZpk2QvFn%2BgwqhYdol%2B%2FF%2F2UiFXBWwtFoAo5G4hY7qGEcLlncpi1PXUTqOuq6DgrvsyQ%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2009%3A00%3A05%20GMT&X-Amz-Signature=fda2f7c6b24d91de30405190349e4c466d9255a4e5f151de2cc212a692f44f30"
$ws.Range("K3").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1750182961683/iBjUYEdqdu-camera-video-segment-driver-1750182964183.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSF4KUIPEM%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T010005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDWqpT3vXbTcJEykZGM%2F7X%2B8cCmNzx24gLc6N7vIEBARAIhALOJOs9Gn6j%2BXidxNTNrq7Bxn%2F5siLlT74MJ6mGD15mIKuYDCJj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgwwGDdKhFzPvUQY3wAqugPXUgcgeUu4qcQIkG3N7shPvpQW7j5wHUCsgW7V53Olq6t0LfSqJfNHT26wH14Bq12VcntUBVt5E2qrV2PBzUgFLpWipERqDdCJeoofeWt1XFjaXBLNYHgUhXAxcTvrMJPf5C5cycHyoE3KORYyzBwlCzNCwbexzk%2BM0iWoauNsOfRGqc9aElrjUtEFZcWh%2Fd5EWgsWXH3SyJQa0C%2F5GLaXRiQgxaS20fMxuwiLnKuD7SJilwLCLmHe7msZKNB2W3n4JXM3jTwbpEzya1PDZhxlXOvuwA%2BkMuZGC9i7PgdlzplqTSkk8sXzVFRtpZA4L1NCItMSBy8AzE5vTnhZJi6wpfOtswcg11UuBe2jdtSWO7qjXTWBHztMrWTVYdk8ghl9dE5vVgPdeeEAkW%2BI5OsOccb5YikxLxgZV1PWFhvdNjl%2FLtaIDTx%2FrG%2FlP9S%2Fwg9kBS1CF%2BzpXfQyPUYtvH8vpFNulkOMACWSQMlHDpAU%2B0sTJifj2plwRRLJYqsx96E8FP2U0ihLZpQC4up%2F1KHOL3EvCVQxbMCma5hOGQhWomUr%2BvfBhCCXgL%2F025TO5cBe4OVu9WkJZfrJMLD8zMIGOqQB7pW6hWXA2tomga2LqDXXQey05yP81WXFX1Yep26oiPmD5hCVc1H1A80szBa4GjLlmAYvMOMvJ3hcZHpy2vIj9KO70quMcwovPnCeBKaoO5DROHvq9oO3NRPMetMR9u5qmXrlQubOcL04T5Ez6O6k8QjzZpk2QvFn%2BgwqhYdol%2B%2FF%2F2UiFXBWwtFoAo5G4hY7qGEcLlncpi1PXUTqOuq6DgrvsyQ%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2009%3A00%3A05%20GMT&X-Amz-Signature=1ff3e05b5c2415b1ea75d0d1501f8eff94229ddf6648286e1314e07ac1d827ad"
$ws.Range("L3").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1750182961683/Foph7ZSROm-camera-video-segment-1750182964183.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSF4KUIPEM%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T010005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDWqpT3vXbTcJEykZGM%2F7X%2B8cCmNzx24gLc6N7vIEBARAIhALOJOs9Gn6j%2BXidxNTNrq7Bxn%2F5siLlT74MJ6mGD15mIKuYDCJj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgwwGDdKhFzPvUQY3wAqugPXUgcgeUu4qcQIkG3N7shPvpQW7j5wHUCsgW7V53Olq6t0LfSqJfNHT26wH14Bq12VcntUBVt5E2qrV2PBzUgFLpWipERqDdCJeoofeWt1XFjaXBLNYHgUhXAxcTvrMJPf5C5cycHyoE3KORYyzBwlCzNCwbexzk%2BM0iWoauNsOfRGqc9aElrjUtEFZcWh%2Fd5EWgsWXH3SyJQa0C%2F5GLaXRiQgxaS20fMxuwiLnKuD7SJilwLCLmHe7msZKNB2W3n4JXM3jTwbpEzya1PDZhxlXOvuwA%2BkMuZGC9i7PgdlzplqTSkk8sXzVFRtpZA4L1NCItMSBy8AzE5vTnhZJi6wpfOtswcg11UuBe2jdtSWO7qjXTWBHztMrWTVYdk8ghl9dE5vVgPdeeEAkW%2BI5OsOccb5YikxLxgZV1PWFhvdNjl%2FLtaIDTx%2FrG%2FlP9S%2Fwg9kBS1CF%2BzpXfQyPUYtvH8vpFNulkOMACWSQMlHDpAU%2B0sTJifj2plwRRLJYqsx96E8FP2U0ihLZpQC4up%2F1KHOL3EvCVQxbMCma5hOGQhWomUr%2BvfBhCCXgL%2F025TO5cBe4OVu9WkJZfrJMLD8zMIGOqQB7pW6hWXA2tomga2LqDXXQey05yP81WXFX1Yep26oiPmD5hCVc1H1A80szBa4GjLlmAYvMOMvJ3hcZHpy2vIj9KO70quMcwovPnCeBKaoO5DROHvq9oO3NRPMetMR9u5qmXrlQubOcL04T5Ez6O6k8QjzZpk2QvFn%2BgwqhYdol%2B%2FF%2F2UiFXBWwtFoAo5G4hY7qGEcLlncpi1PXUTqOuq6DgrvsyQ%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2009%3A00%3A05%20GMT&X-Amz-Signature=512c99290a31cf0315a5067742f24a9fad759dd6eef5e4934cf411305a008cd5"
